# Auto-generated Excel COM-interop edit script
# Applies the weekly CompStat data refresh (new crime data collected)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/issue number and week-covering date range) ---
$ws.Range("A8").Value = "Volume 31   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/8/2024  Through  4/14/2024"

# --- Updated weekly crime statistics table (rows 15-31) ---
# Row 15
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("L15").Value = -33.333333333333
$ws.Range("L15").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 16
$ws.Range("C16").Value = 2
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("D16").Value = 5
$ws.Range("D16").NumberFormat = '#,##0'
$ws.Range("E16").Value = -60
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F16").Value = 6
$ws.Range("F16").NumberFormat = '#,##0'
$ws.Range("G16").Value = 9
$ws.Range("G16").NumberFormat = '#,##0'
$ws.Range("H16").Value = -33.333333333333
$ws.Range("H16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I16").Value = 32
$ws.Range("I16").NumberFormat = '#,##0'
$ws.Range("J16").Value = 36
$ws.Range("J16").NumberFormat = '#,##0'
$ws.Range("K16").Value = -11.111111111111
$ws.Range("K16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L16").Value = -23.809523809523
$ws.Range("L16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M16").Value = 60
$ws.Range("M16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N16").Value = -88.059701492537
$ws.Range("N16").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 17
$ws.Range("C17").Value = 3
$ws.Range("C17").NumberFormat = '#,##0'
$ws.Range("D17").Value = 4
$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("E17").Value = -25
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F17").Value = 18
$ws.Range("F17").NumberFormat = '#,##0'
$ws.Range("G17").Value = 9
$ws.Range("G17").NumberFormat = '#,##0'
$ws.Range("H17").Value = 100
$ws.Range("H17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I17").Value = 45
$ws.Range("I17").NumberFormat = '#,##0'
$ws.Range("J17").Value = 36
$ws.Range("J17").NumberFormat = '#,##0'
$ws.Range("K17").Value = 25
$ws.Range("K17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L17").Value = 36.363636363636
$ws.Range("L17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M17").Value = 221.428571428571
$ws.Range("M17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N17").Value = -13.461538461538
$ws.Range("N17").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("C18").NumberFormat = '#,##0'
$ws.Range("D18").Value = 2
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("E18").Value = 50
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G18").Value = 16
$ws.Range("G18").NumberFormat = '#,##0'
$ws.Range("H18").Value = -31.25
$ws.Range("H18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I18").Value = 62
$ws.Range("I18").NumberFormat = '#,##0'
$ws.Range("J18").Value = 48
$ws.Range("J18").NumberFormat = '#,##0'
$ws.Range("K18").Value = 29.166666666666
$ws.Range("K18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L18").Value = -26.190476190476
$ws.Range("L18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M18").Value = 8.771929824561
$ws.Range("M18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N18").Value = -75
$ws.Range("N18").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 19
$ws.Range("C19").Value = 21
$ws.Range("C19").NumberFormat = '#,##0'
$ws.Range("D19").Value = 31
$ws.Range("D19").NumberFormat = '#,##0'
$ws.Range("E19").Value = -32.258064516129
$ws.Range("E19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F19").Value = 84
$ws.Range("F19").NumberFormat = '#,##0'
$ws.Range("G19").Value = 103
$ws.Range("G19").NumberFormat = '#,##0'
$ws.Range("H19").Value = -18.446601941747
$ws.Range("H19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I19").Value = 299
$ws.Range("I19").NumberFormat = '#,##0'
$ws.Range("J19").Value = 329
$ws.Range("J19").NumberFormat = '#,##0'
$ws.Range("K19").Value = -9.118541033434
$ws.Range("K19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L19").Value = -5.379746835443
$ws.Range("L19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M19").Value = -3.236245954692
$ws.Range("M19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N19").Value = -68.592436974789
$ws.Range("N19").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 20
$ws.Range("C20").Value = 2
$ws.Range("C20").NumberFormat = '#,##0'
$ws.Range("D20").Value = 4
$ws.Range("D20").NumberFormat = '#,##0'
$ws.Range("E20").Value = -50
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F20").Value = 6
$ws.Range("F20").NumberFormat = '#,##0'
$ws.Range("G20").Value = 10
$ws.Range("G20").NumberFormat = '#,##0'
$ws.Range("H20").Value = -40
$ws.Range("H20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I20").Value = 10
$ws.Range("I20").NumberFormat = '#,##0'
$ws.Range("J20").Value = 17
$ws.Range("J20").NumberFormat = '#,##0'
$ws.Range("K20").Value = -41.176470588235
$ws.Range("K20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L20").Value = 11.111111111111
$ws.Range("L20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M20").Value = 150
$ws.Range("M20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N20").Value = -95.850622406639
$ws.Range("N20").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 21
$ws.Range("C21").Value = 31
$ws.Range("C21").NumberFormat = '#,##0'
$ws.Range("D21").Value = 46
$ws.Range("D21").NumberFormat = '#,##0'
$ws.Range("E21").Value = -32.608695652173
$ws.Range("E21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("G21").Value = 150
$ws.Range("G21").NumberFormat = '#,##0'
$ws.Range("H21").Value = -16.666666666666
$ws.Range("H21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("I21").Value = 453
$ws.Range("I21").NumberFormat = '#,##0'
$ws.Range("J21").Value = 471
$ws.Range("J21").NumberFormat = '#,##0'
$ws.Range("K21").Value = -3.821656050955
$ws.Range("K21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("L21").Value = -7.551020408163
$ws.Range("L21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("M21").Value = 12.128712871287
$ws.Range("M21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("N21").Value = -74.406779661017
$ws.Range("N21").NumberFormat = '#,##0.00;"-"#,##0.00'
# Row 22
$ws.Range("C22").Value = 3
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 3
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F22").Value = 6
$ws.Range("F22").NumberFormat = '#,##0'
$ws.Range("G22").Value = 4
$ws.Range("G22").NumberFormat = '#,##0'
$ws.Range("H22").Value = 50
$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I22").Value = 20
$ws.Range("I22").NumberFormat = '#,##0'
$ws.Range("J22").Value = 25
$ws.Range("J22").NumberFormat = '#,##0'
$ws.Range("K22").Value = -20
$ws.Range("K22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L22").Value = -35.483870967741
$ws.Range("L22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M22").Value = -13.043478260869
$ws.Range("M22").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 24
$ws.Range("C24").Value = 98
$ws.Range("C24").NumberFormat = '#,##0'
$ws.Range("D24").Value = 66
$ws.Range("D24").NumberFormat = '#,##0'
$ws.Range("E24").Value = 48.484848484848
$ws.Range("E24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F24").Value = 328
$ws.Range("F24").NumberFormat = '#,##0'
$ws.Range("G24").Value = 263
$ws.Range("G24").NumberFormat = '#,##0'
$ws.Range("H24").Value = 24.714828897338
$ws.Range("H24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I24").Value = 1173
$ws.Range("I24").NumberFormat = '#,##0'
$ws.Range("J24").Value = 1063
$ws.Range("J24").NumberFormat = '#,##0'
$ws.Range("K24").Value = 10.348071495766
$ws.Range("K24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L24").Value = 13.994169096209
$ws.Range("L24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M24").Value = 165.384615384615
$ws.Range("M24").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 25
$ws.Range("C25").Value = 92
$ws.Range("C25").NumberFormat = '#,##0'
$ws.Range("D25").Value = 72
$ws.Range("D25").NumberFormat = '#,##0'
$ws.Range("E25").Value = 27.777777777777
$ws.Range("E25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F25").Value = 316
$ws.Range("F25").NumberFormat = '#,##0'
$ws.Range("G25").Value = 273
$ws.Range("G25").NumberFormat = '#,##0'
$ws.Range("H25").Value = 15.750915750915
$ws.Range("H25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I25").Value = 1158
$ws.Range("I25").NumberFormat = '#,##0'
$ws.Range("J25").Value = 1087
$ws.Range("J25").NumberFormat = '#,##0'
$ws.Range("K25").Value = 6.531738730450
$ws.Range("K25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L25").Value = 12.427184466019
$ws.Range("L25").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 26
$ws.Range("D26").Value = 7
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("E26").Value = 28.571428571428
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F26").Value = 38
$ws.Range("F26").NumberFormat = '#,##0'
$ws.Range("G26").Value = 22
$ws.Range("G26").NumberFormat = '#,##0'
$ws.Range("H26").Value = 72.727272727272
$ws.Range("H26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I26").Value = 123
$ws.Range("I26").NumberFormat = '#,##0'
$ws.Range("J26").Value = 86
$ws.Range("J26").NumberFormat = '#,##0'
$ws.Range("K26").Value = 43.023255813953
$ws.Range("K26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L26").Value = 41.379310344827
$ws.Range("L26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M26").Value = 83.582089552238
$ws.Range("M26").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = '#,##0'
$ws.Range("G27").Value = 4
$ws.Range("G27").NumberFormat = '#,##0'
$ws.Range("H27").Value = -75
$ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I27").Value = 7
$ws.Range("I27").NumberFormat = '#,##0'
$ws.Range("J27").Value = 6
$ws.Range("J27").NumberFormat = '#,##0'
$ws.Range("K27").Value = 16.666666666666
$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 28
$ws.Range("C28").Value = 2
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 5
$ws.Range("F28").NumberFormat = '#,##0'
$ws.Range("H28").Value = 25
$ws.Range("H28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I28").Value = 28
$ws.Range("I28").NumberFormat = '#,##0'
$ws.Range("K28").Value = 40
$ws.Range("K28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L28").Value = 21.739130434782
$ws.Range("L28").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 31
$ws.Range("C31").Value = 1
$ws.Range("C31").NumberFormat = '#,##0'
$ws.Range("F31").Value = 3
$ws.Range("F31").NumberFormat = '#,##0'
$ws.Range("I31").Value = 4
$ws.Range("I31").NumberFormat = '#,##0'
$ws.Range("K31").Value = -20
$ws.Range("K31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L31").Value = -60
$ws.Range("L31").NumberFormat = '#,##0.0;"-"#,##0.0'
